# Apply updated Leve profit/price figures per the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1523.5714
$ws.Range("I2").Value = 355
$ws.Range("K2").Value = 355
$ws.Range("M2").Value = -242
$ws.Range("H19").Value = 936.3077
$ws.Range("J19").Value = 750.125
$ws.Range("L19").Value = 750.125
$ws.Range("N19").Value = -1100.125
$ws.Range("H107").Value = 653.5714
$ws.Range("I107").Value = 710.5454999999999
$ws.Range("J107").Value = 444.66666
$ws.Range("K107").Value = 710.5454999999999
$ws.Range("L107").Value = 444.66666
$ws.Range("M107").Value = 1209.4545
$ws.Range("N107").Value = -4284.66666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1657.091
$ws.Range("I2").Value = 921.3
$ws.Range("K2").Value = 921.3
$ws.Range("M2").Value = -808.3
$ws.Range("H32").Value = 1039.5625
$ws.Range("I32").Value = 1039.5625
$ws.Range("K32").Value = 1039.5625
$ws.Range("M32").Value = -752.5625
$ws.Range("H45").Value = 2222.875
$ws.Range("I45").Value = 1207.4445
$ws.Range("J45").Value = 3528.4285
$ws.Range("K45").Value = 1207.4445
$ws.Range("L45").Value = 3528.4285
$ws.Range("M45").Value = -830.4445000000001
$ws.Range("N45").Value = -4282.4285
$ws.Range("H74").Value = 6561.5
$ws.Range("I74").Value = 6584.857
$ws.Range("K74").Value = 6584.857
$ws.Range("M74").Value = -5710.857
$ws.Range("H77").Value = 6561.5
$ws.Range("I77").Value = 6584.857
$ws.Range("K77").Value = 32924.285
$ws.Range("M77").Value = -28556.285
$ws.Range("H97").Value = 700.4
$ws.Range("I97").Value = 834.8333
$ws.Range("K97").Value = 834.8333
$ws.Range("M97").Value = -338.8333
$ws.Range("H110").Value = 1833
$ws.Range("I110").Value = 1499.5
$ws.Range("J110").Value = 2500
$ws.Range("K110").Value = 1499.5
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = 545.5
$ws.Range("N110").Value = -6590
$ws.Range("H116").Value = 1657.091
$ws.Range("I116").Value = 921.3
$ws.Range("K116").Value = 921.3
$ws.Range("M116").Value = 1372.7
$ws.Range("H122").Value = 983.1667
$ws.Range("I122").Value = 983.1667
$ws.Range("K122").Value = 2949.5001
$ws.Range("M122").Value = -499.5001000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1657.091
$ws.Range("I3").Value = 921.3
$ws.Range("K3").Value = 921.3
$ws.Range("M3").Value = -807.3
$ws.Range("H105").Value = 1676.8334
$ws.Range("I105").Value = 1676.8334
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1676.8334
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 70.16660000000002
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 7439
$ws.Range("I107").Value = 5975
$ws.Range("K107").Value = 5975
$ws.Range("M107").Value = -4055
$ws.Range("H122").Value = 59999
$ws.Range("J122").Value = 59999
$ws.Range("L122").Value = 59999
$ws.Range("N122").Value = -69799
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3833
$ws.Range("I86").Value = 3749.5
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 3749.5
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -2626.5
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 3833
$ws.Range("I89").Value = 3749.5
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 18747.5
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -13131.5
$ws.Range("N89").Value = -31232
$ws.Range("H122").Value = 969.5
$ws.Range("I122").Value = 608
$ws.Range("K122").Value = 1824
$ws.Range("M122").Value = 626
$ws.Range("H134").Value = 3250.375
$ws.Range("I134").Value = 3250.375
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9751.125
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7216.125
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1230.6
$ws.Range("J107").Value = 1322.8889
$ws.Range("L107").Value = 3968.6667
$ws.Range("N107").Value = -7808.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 1631.75
$ws.Range("I36").Value = 2355.6667
$ws.Range("J36").Value = 1197.4
$ws.Range("K36").Value = 2355.6667
$ws.Range("L36").Value = 1197.4
$ws.Range("M36").Value = -1870.6667
$ws.Range("N36").Value = -2167.4
$ws.Range("H97").Value = 733.8889
$ws.Range("J97").Value = 924.75
$ws.Range("L97").Value = 924.75
$ws.Range("N97").Value = -1916.75
$ws.Range("H102").Value = 4989.2
$ws.Range("I102").Value = 5236.5
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 5236.5
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -3614.5
$ws.Range("N102").Value = -7244
$ws.Range("H113").Value = 7261.3335
$ws.Range("J113").Value = 7261.3335
$ws.Range("L113").Value = 7261.3335
$ws.Range("N113").Value = -11601.3335
$ws.Range("H122").Value = 3010.5
$ws.Range("I122").Value = 3010.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9031.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6581.5
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 3999.5
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 4999
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 4999
$ws.Range("M33").Value = -2710
$ws.Range("N33").Value = -5579
$ws.Range("H35").Value = 1795.3
$ws.Range("I35").Value = 1324.2222
$ws.Range("J35").Value = 6035
$ws.Range("K35").Value = 1324.2222
$ws.Range("L35").Value = 6035
$ws.Range("M35").Value = -988.2221999999999
$ws.Range("N35").Value = -6707
$ws.Range("H40").Value = 4187.6665
$ws.Range("I40").Value = 4187.6665
$ws.Range("K40").Value = 4187.6665
$ws.Range("M40").Value = -4051.6665
$ws.Range("H61").Value = 4523.5
$ws.Range("I61").Value = 3372.5
$ws.Range("J61").Value = 6250
$ws.Range("K61").Value = 3372.5
$ws.Range("L61").Value = 6250
$ws.Range("M61").Value = -3170.5
$ws.Range("N61").Value = -6654
$ws.Range("H82").Value = 2857.5908
$ws.Range("I82").Value = 1484.8462
$ws.Range("J82").Value = 4840.4443
$ws.Range("K82").Value = 1484.8462
$ws.Range("L82").Value = 4840.4443
$ws.Range("M82").Value = -1123.8462
$ws.Range("N82").Value = -5562.4443
$ws.Range("H85").Value = 2857.5908
$ws.Range("I85").Value = 1484.8462
$ws.Range("J85").Value = 4840.4443
$ws.Range("K85").Value = 1484.8462
$ws.Range("L85").Value = 4840.4443
$ws.Range("M85").Value = -236.8462
$ws.Range("N85").Value = -7336.4443
$ws.Range("H100").Value = 9142.857
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459
$ws.Range("H113").Value = 4523.5
$ws.Range("I113").Value = 3372.5
$ws.Range("J113").Value = 6250
$ws.Range("K113").Value = 3372.5
$ws.Range("L113").Value = 6250
$ws.Range("M113").Value = -1202.5
$ws.Range("N113").Value = -10590
$ws.Range("H132").Value = 3131.75
$ws.Range("I132").Value = 3091
$ws.Range("K132").Value = 9273
$ws.Range("M132").Value = -6743

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 937.86664
$ws.Range("I100").Value = 706.1818
$ws.Range("K100").Value = 1412.3636
$ws.Range("M100").Value = -871.3635999999999
$ws.Range("H132").Value = 956
$ws.Range("I132").Value = 525.44446
$ws.Range("K132").Value = 1576.33338
$ws.Range("M132").Value = 953.66662
